$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6974
$ws1.Range("F4").Value = 456
$ws1.Range("F6").Value = 550
$ws1.Range("F7").Value = 128
$ws1.Range("F10").Value = 12
$ws1.Range("F12").Value = 196
$ws1.Range("F14").Value = 25
$ws1.Range("F17").Value = 3584
$ws1.Range("F20").Value = 81
$ws1.Range("F23").Value = 2201
$ws1.Range("F25").Value = 233
$ws1.Range("F32").Value = 119
$ws1.Range("F33").Value = 58

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6974
$ws4.Range("F4").Value = 456
$ws4.Range("F7").Value = 550
$ws4.Range("F8").Value = 128
$ws4.Range("F11").Value = 12
$ws4.Range("F13").Value = 196
$ws4.Range("F15").Value = 25
$ws4.Range("F18").Value = 3584
$ws4.Range("F21").Value = 81
$ws4.Range("F24").Value = 2201
$ws4.Range("F26").Value = 233
$ws4.Range("F33").Value = 119
$ws4.Range("F34").Value = 58

$wb.Save()
